# Continuing the development of the Monty JSON Schema
# Add the missing "haztype*" sub-category labels in column C for the rows
# that previously only had a hazard name in column A (Civil Unrest,
# Transport Accident, Chemical Emergency, Insect Infestation and
# Biological Emergency).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "haztypesoc"
$ws.Range("C19").Value = "haztypetech"
$ws.Range("C21").Value = "haztypechem"
$ws.Range("C22").Value = "haztypebio"
$ws.Range("C23").Value = "haztypebio"

# Rows whose content changed/grew get a recalculated (slightly tighter)
# row height in the source workbook.
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(19).RowHeight = 13.8
$ws.Rows.Item(21).RowHeight = 13.8
$ws.Rows.Item(22).RowHeight = 13.8
$ws.Rows.Item(23).RowHeight = 13.8

# The author's cursor ended up on the newly-edited C14 cell.
$ws.Range("C14").Select()
